$d = $word.ActiveDocument

function Split-At($pos) {
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add("ZZSPLITTMP", $r) | Out-Null
    $d.Bookmarks.Item("ZZSPLITTMP").Delete()
}

# ------------------------------------------------------------------
# 1. Locate the run to edit: " experience working with Player Input,
#    animations and AI on Action Games. "
# ------------------------------------------------------------------
$findRng = $d.Content
$find = $findRng.Find
$find.ClearFormatting()
$find.Execute("working", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$workingStart = $findRng.Start
$workingEnd = $findRng.End   # end of the word "working" (not including trailing space)

# Pre-split the run so that "working " becomes an isolated run before we cut it
Split-At $workingStart
Split-At ($workingEnd + 1)   # +1 to include the trailing space in the cut piece

# ------------------------------------------------------------------
# 2. Cut "working " out and paste it back in immediately before
#    "on Action Games."
# ------------------------------------------------------------------
$cutRng = $d.Range($workingStart, $workingEnd + 1)
$cutRng.Cut() | Out-Null

$findRng2 = $d.Content
$find2 = $findRng2.Find
$find2.ClearFormatting()
$find2.Execute("on Action Games.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$onActionStart = $findRng2.Start

$pasteRng = $d.Range($onActionStart, $onActionStart)
$pasteRng.Paste() | Out-Null

# ------------------------------------------------------------------
# 3. The Paste operation merges everything from "on Action Games."
#    through the end of the paragraph into a single run. Re-split it
#    back into the original run boundaries.
# ------------------------------------------------------------------
$findRng3 = $d.Content
$find3 = $findRng3.Find
$find3.ClearFormatting()
$find3.Execute("working on Action Games.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base = $findRng3.Start   # start of "working "

# boundaries, relative to $base (start of "working "):
#    0 -> "working "                              (already isolated above)
#    8 -> "on Action Games. "
#   25 -> "I am currently a graduate student..."
#  112 -> "for"
#  115 -> " my "
#  119 -> "EAE "
#  123 -> "- "
#  125 -> "Game Engineering"
#  141 -> " Master's"
#  150 -> " degree. And I am "
#  168 -> "beyond "
#  175 -> "thrilled"
#  183 -> " to be applying for the Gameplay Programm"
#  224 -> "er"
#  226 -> " Intern position at Santa Monica Studio."
$boundaries = @(8, 25, 112, 115, 119, 123, 125, 141, 150, 168, 175, 183, 224, 226)

foreach ($b in $boundaries) {
    Split-At ($base + $b)
}

# ------------------------------------------------------------------
# 4. Move the _GoBack bookmark to sit between "working " and
#    "on Action Games. "
# ------------------------------------------------------------------
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

$bmPos = $base + 8
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null

Write-Host "Done"
